# This workbook lists public-exposure-site case alerts, one site per row
# (columns: Location, Site, Exposure period, Notes). The update:
#   - removes six sites that are no longer current exposure sites
#     (Clayton South; Heatherton; Moorabbin Airport; and all three
#     Springvale entries), and
#   - corrects the exposure-period end time for the Pascoe Vale /
#     "Elite Swimming Pascoe Vale" entry from 6:00pm to 5:30pm.
# Because rows below a deleted row shift up, the row numbers below are
# expressed against the *original* layout and deletions are applied from
# the bottom up so earlier row indices stay valid while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Pascoe Vale / Elite Swimming exposure period (row 24, column C)
# before any rows are removed, while row 24 still refers to that record.
$ws.Cells.Item(24, 3).Value = "5pm - 5:30pm 8/2/2021"

# Delete the obsolete rows, highest row number first so that the indices
# of rows still to be deleted are unaffected by earlier deletions.
$ws.Rows(28).Delete()   # Springvale - Woolworths Springvale
$ws.Rows(27).Delete()   # Springvale - Sharetea Springvale
$ws.Rows(26).Delete()   # Springvale - Bunnings Springvale
$ws.Rows(23).Delete()   # Moorabbin Airport - Lululemon, DFO Moorabbin
$ws.Rows(13).Delete()   # Heatherton - Melbourne Golf Academy
$ws.Rows(8).Delete()    # Clayton South - Nakama Workshop
